$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.962.81"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.486.52"
$ws.Range("E3").Value = "  -1.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.57"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.90"
$ws.Range("E6").Value = "  -3.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.517"
$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -5.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.52"
$ws.Range("E10").Value = "  -4.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0798"
$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.02"
$ws.Range("E14").Value = "  -3.70%  "

$ws.Range("D15").Value = "2.881.27"
$ws.Range("E15").Value = "  -1.38%  "

$ws.Range("D16").Value = "2.492.99"
$ws.Range("E16").Value = "  -1.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.825"
$ws.Range("E17").Value = "  -4.46%  "

$ws.Range("D18").Value = "47.866.19"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.95"
$ws.Range("E19").Value = "  +8.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").Value = "0.0₃0926"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.85"
$ws.Range("E23").Value = "  -2.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.76"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  -3.33%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.62"
$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  +4.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  -5.11%  "

$ws.Range("E30").Value = "  -6.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.25"
$ws.Range("E31").Value = "  -3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.14"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.25"
$ws.Range("E34").Value = "  -2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.82"
$ws.Range("E35").Value = "  -5.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0769"
$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.55"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  -5.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.18"
$ws.Range("E40").Value = "  +3.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.24"
$ws.Range("E41").Value = "  -1.71%  "

$ws.Range("E42").Value = "  -2.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.21"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").Value = "1.994.01"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.10"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.85"
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.12"
$ws.Range("E50").Value = "  -2.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.28"
$ws.Range("E51").Value = "  -2.69%  "
